$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 34/35 swap: LidoDAOToken <-> ImmutableX (ranking positions exchanged) ---
$ws.Range("B34").Value = "'ImmutableX"
$ws.Range("C34").Value = "'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D34").Value = "'0.7321"
$ws.Range("E34").Value = "'  -1.41%  "

$ws.Range("B35").Value = "'LidoDAOToken"
$ws.Range("C35").Value = "'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D35").Value = "'1.812"
$ws.Range("E35").Value = "'  -1.58%  "

# --- Regular per-row Price (D) / Volume(1h) (E) updates ---
$ws.Range("D2").Value = "'29.069.72"
$ws.Range("E2").Value = "'  -0.10%  "
$ws.Range("D3").Value = "'1.821.02"
$ws.Range("E3").Value = "'  -0.73%  "
$ws.Range("D4").Value = "'1.000"
$ws.Range("E4").Value = "'  -0.18%  "
$ws.Range("D5").Value = "'241.37"
$ws.Range("E5").Value = "'  -0.94%  "
$ws.Range("D6").Value = "'0.6136"
$ws.Range("E6").Value = "'  -2.28%  "
$ws.Range("D7").Value = "'1.000"
$ws.Range("E7").Value = "'  -0.14%  "
$ws.Range("D8").Value = "'0.07314"
$ws.Range("E8").Value = "'  -2.63%  "
$ws.Range("D9").Value = "'0.2879"
$ws.Range("E9").Value = "'  -1.45%  "
$ws.Range("D10").Value = "'22.82"
$ws.Range("E10").Value = "'  -1.65%  "
$ws.Range("D11").Value = "'0.07650"
$ws.Range("E11").Value = "'  -0.36%  "
$ws.Range("D12").Value = "'1.815.51"
$ws.Range("E12").Value = "'  -0.99%  "
$ws.Range("D13").Value = "'4.936"
$ws.Range("E13").Value = "'  -1.39%  "
$ws.Range("D14").Value = "'0.6590"
$ws.Range("E14").Value = "'  -1.23%  "
$ws.Range("D15").Value = "'81.41"
$ws.Range("E15").Value = "'  -1.55%  "
$ws.Range("D16").Value = "'0.000008941"
$ws.Range("E16").Value = "'  -4.65%  "
$ws.Range("D17").Value = "'5.832"
$ws.Range("E17").Value = "'  -2.54%  "
$ws.Range("D18").Value = "'29.049.02"
$ws.Range("E18").Value = "'  -0.14%  "
$ws.Range("D19").Value = "'2.063.41"
$ws.Range("E19").Value = "'  -0.78%  "
$ws.Range("D20").Value = "'237.71"
$ws.Range("E20").Value = "'  +6.51%  "
$ws.Range("D21").Value = "'12.41"
$ws.Range("E21").Value = "'  -1.37%  "
$ws.Range("D22").Value = "'1.0000"
$ws.Range("E22").Value = "'  -0.37%  "
$ws.Range("D23").Value = "'7.113"
$ws.Range("E23").Value = "'  +0.20%  "
$ws.Range("E24").Value = "'  -0.22%  "
$ws.Range("D25").Value = "'157.58"
$ws.Range("E25").Value = "'  -1.28%  "
$ws.Range("D26").Value = "'0.1404"
$ws.Range("E26").Value = "'  +0.39%  "
$ws.Range("D27").Value = "'8.404"
$ws.Range("E27").Value = "'  -0.99%  "
$ws.Range("D28").Value = "'17.58"
$ws.Range("E28").Value = "'  -1.82%  "
$ws.Range("E29").Value = "'  -1.40%  "
$ws.Range("D30").Value = "'0.05568"
$ws.Range("E30").Value = "'  -1.75%  "
$ws.Range("D31").Value = "'4.083"
$ws.Range("E31").Value = "'  +0.15%  "
$ws.Range("D32").Value = "'4.094"
$ws.Range("E32").Value = "'  -1.43%  "
$ws.Range("D33").Value = "'1.207"
$ws.Range("E33").Value = "'  +0.17%  "
$ws.Range("D36").Value = "'1.127"
$ws.Range("E36").Value = "'  -1.13%  "
$ws.Range("D37").Value = "'2.617"
$ws.Range("E37").Value = "'  -1.98%  "
$ws.Range("D38").Value = "'2.846"
$ws.Range("E38").Value = "'  +2.92%  "
$ws.Range("D39").Value = "'1.201.90"
$ws.Range("E39").Value = "'  -1.63%  "
$ws.Range("D40").Value = "'0.01755"
$ws.Range("E40").Value = "'  -1.41%  "
$ws.Range("D41").Value = "'6.357"
$ws.Range("E41").Value = "'  -2.67%  "
$ws.Range("D42").Value = "'0.8899"
$ws.Range("E42").Value = "'  -0.39%  "
$ws.Range("E43").Value = "'  -0.09%  "
$ws.Range("D44").Value = "'100.56"
$ws.Range("E44").Value = "'  -1.37%  "
$ws.Range("D45").Value = "'1.966.69"
$ws.Range("E45").Value = "'  -0.63%  "
$ws.Range("D46").Value = "'64.51"
$ws.Range("E46").Value = "'  -2.00%  "
$ws.Range("D49").Value = "'9.026"
$ws.Range("E49").Value = "'  -0.11%  "
$ws.Range("D50").Value = "'0.3978"
$ws.Range("E50").Value = "'  -2.40%  "
$ws.Range("D51").Value = "'0.05785"
$ws.Range("E51").Value = "'  -0.71%  "

# --- Row 47/48 swap: BabyDogeCoin <-> Mantle (ranking positions exchanged) ---
$ws.Range("B47").Value = "'Mantle"
$ws.Range("C47").Value = "'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D47").Value = "'0.5090"
$ws.Range("E47").Value = "'  -0.15%  "

$ws.Range("B48").Value = "'BabyDogeCoin"
$ws.Range("C48").Value = "'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D48").Value = "'0.00000000120"
$ws.Range("E48").Value = "'  -4.85%  "
